$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.447.72"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.447.47"
$ws.Range("E3").Value = "  -2.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.00"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.41"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "2.447.13"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("E10").Value = "  -6.39%  "
$ws.Range("E12").Value = "  -5.89%  "
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "2.892.05"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "68.257.29"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.59"
$ws.Range("E17").Value = "  -5.17%  "
$ws.Range("D18").Value = "2.476.62"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.96"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.98"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.21"
$ws.Range("E21").Value = "  -5.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.79"
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.85"
$ws.Range("E24").Value = "  -6.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.78"
$ws.Range("E25").Value = "  -4.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.75"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("D27").Value = "2.569.52"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.22"
$ws.Range("E29").Value = "  -7.12%  "
$ws.Range("D30").Value = "0.0₃0837"
$ws.Range("E30").Value = "  -6.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  -7.21%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -4.67%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.67"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "426.63"
$ws.Range("E35").Value = "  -7.75%  "
$ws.Range("E36").Value = "  +109.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.17"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.00"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -5.77%  "
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.49"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.53"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.07"
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("E46").Value = "  -7.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.52"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0716"
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.486"
$ws.Range("E50").Value = "  -6.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.558"
$ws.Range("E51").Value = "  -3.42%  "
